$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 duplicated row 2's entry but with a later date (13-FEB-26).
# Merge it into row 2 by copying row 3's date onto row 2 (this preserves
# the original cell formatting/style), then remove row 3 entirely.
$ws.Range("A3").Copy($ws.Range("A2"))
$ws.Rows("3").Delete()
